$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to read a row's varying values (A, B, E, F, G, H, Q, R, Z, AB)
# Note: use Value2 (not Value) for both get and set - the Value getter in
# this runtime does not reliably return the underlying cell value.
function Get-RowData($row) {
    $data = @{}
    $data.A  = $ws.Range("A$row").Value2
    $data.B  = $ws.Range("B$row").Value2
    $data.E  = $ws.Range("E$row").Value2
    $data.F  = $ws.Range("F$row").Value2
    $data.G  = $ws.Range("G$row").Value2
    $data.H  = $ws.Range("H$row").Value2
    $data.Q  = $ws.Range("Q$row").Value2
    $data.R  = $ws.Range("R$row").Value2
    $data.Z  = $ws.Range("Z$row").Value2
    $data.AB = $ws.Range("AB$row").Value2
    return $data
}

function Set-RowData($row, $data) {
    $ws.Range("A$row").Value2  = $data.A
    $ws.Range("B$row").Value2  = $data.B
    $ws.Range("E$row").Value2  = $data.E
    $ws.Range("F$row").Value2  = $data.F
    $ws.Range("G$row").Value2  = $data.G
    $ws.Range("H$row").Value2  = $data.H
    $ws.Range("Q$row").Value2  = $data.Q
    $ws.Range("R$row").Value2  = $data.R
    $ws.Range("Z$row").Value2  = $data.Z
    $ws.Range("AB$row").Value2 = $data.AB
}

# Row 3 <-> Row 4 (full swap of the varying identifying columns)
$row3 = Get-RowData 3
$row4 = Get-RowData 4
Set-RowData 3 $row4
Set-RowData 4 $row3

# Row 8 <-> Row 9 (full swap)
$row8 = Get-RowData 8
$row9 = Get-RowData 9
Set-RowData 8 $row9
Set-RowData 9 $row8

# Row 13 <-> Row 14 (full swap)
$row13 = Get-RowData 13
$row14 = Get-RowData 14
Set-RowData 13 $row14
Set-RowData 14 $row13

# Row 17 -> 18 -> 19 -> 17 (3-way rotation)
$row17 = Get-RowData 17
$row18 = Get-RowData 18
$row19 = Get-RowData 19
Set-RowData 18 $row17
Set-RowData 19 $row18
Set-RowData 17 $row19
